$wb = $excel.ActiveWorkbook

# --- Sheet "general" ---
$ws = $wb.Worksheets.Item("general")
$ws.Range("G6").Value = 3.4
$ws.Range("G7").Select()

# --- Sheet "window_size" ---
$ws = $wb.Worksheets.Item("window_size")
$ws.Range("E4").Value = 3.4
$ws.Range("E6").Value = 3.4
$ws.Range("E7").Value = 3.4
$ws.Range("E8").Value = 3.4
$ws.Range("E8").Select()

# --- Sheet "scf_size" ---
$ws = $wb.Worksheets.Item("scf_size")
$ws.Range("G4").Value = 3.3
$ws.Range("G6").Value = 3.4
$ws.Range("G7").Value = 3.4
$ws.Range("G8").Value = 3.4
$ws.Range("G8").Select()

# --- Sheet "w2v_size" ---
$ws = $wb.Worksheets.Item("w2v_size")
$ws.Range("F7").Value = 3.4
$ws.Range("G9").Value = 7.6
$ws.Range("G10").Select()

# --- Restore "general" as the active sheet (it was the tab selected originally) ---
$ws = $wb.Worksheets.Item("general")
$ws.Activate()
$ws.Range("G7").Select()
